$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 46600
$ws.Range("J3").Value = 46600
$ws.Range("L3").Value = 46600
$ws.Range("N3").Value = -46828
$ws.Range("H33").Value = 212.53572
$ws.Range("J33").Value = 126
$ws.Range("L33").Value = 126
$ws.Range("N33").Value = -584
$ws.Range("H100").Value = 2006.9333
$ws.Range("I100").Value = 1262.5
$ws.Range("J100").Value = 2857.7144
$ws.Range("K100").Value = 1262.5
$ws.Range("L100").Value = 2857.7144
$ws.Range("M100").Value = -721.5
$ws.Range("N100").Value = -3939.7144
$ws.Range("H102").Value = 46600
$ws.Range("J102").Value = 46600
$ws.Range("L102").Value = 46600
$ws.Range("N102").Value = -53090
$ws.Range("H105").Value = 20960
$ws.Range("J105").Value = 20960
$ws.Range("L105").Value = 20960
$ws.Range("N105").Value = -27948
$ws.Range("H136").Value = 44917
$ws.Range("J136").Value = 44917
$ws.Range("L136").Value = 44917
$ws.Range("N136").Value = -55117
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1809.4117
$ws.Range("I61").Value = 1521
$ws.Range("J61").Value = 3155.3333
$ws.Range("K61").Value = 1521
$ws.Range("L61").Value = 3155.3333
$ws.Range("M61").Value = -1309
$ws.Range("N61").Value = -3579.3333
$ws.Range("H74").Value = 1238.7273
$ws.Range("J74").Value = 1528.4286
$ws.Range("L74").Value = 1528.4286
$ws.Range("N74").Value = -3276.4286
$ws.Range("H77").Value = 1238.7273
$ws.Range("J77").Value = 1528.4286
$ws.Range("L77").Value = 7642.143
$ws.Range("N77").Value = -16378.143
$ws.Range("H132").Value = 1054612.6
$ws.Range("I132").Value = 2000722.8
$ws.Range("K132").Value = 6002168.4
$ws.Range("M132").Value = -5999638.4
$ws.Range("H133").Value = 41566
$ws.Range("J133").Value = 41566
$ws.Range("L133").Value = 41566
$ws.Range("N133").Value = -46626
$ws.Range("H134").Value = 64762.5
$ws.Range("J134").Value = 64762.5
$ws.Range("L134").Value = 64762.5
$ws.Range("N134").Value = -74902.5
$ws.Range("H136").Value = 1809.4117
$ws.Range("I136").Value = 1521
$ws.Range("J136").Value = 3155.3333
$ws.Range("K136").Value = 4563
$ws.Range("L136").Value = 9465.999899999999
$ws.Range("M136").Value = -2013
$ws.Range("N136").Value = -14565.9999
$ws.Range("H139").Value = 85689.37
$ws.Range("J139").Value = 85689.37
$ws.Range("L139").Value = 85689.37
$ws.Range("N139").Value = -95969.37

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1146.9
$ws.Range("I94").Value = 951.8261
$ws.Range("K94").Value = 951.8261
$ws.Range("M94").Value = -500.8261
$ws.Range("H134").Value = 628090.6
$ws.Range("I134").Value = 898040.25
$ws.Range("J134").Value = 3832.125
$ws.Range("K134").Value = 2694120.75
$ws.Range("L134").Value = 11496.375
$ws.Range("M134").Value = -2691585.75
$ws.Range("N134").Value = -16566.375
$ws.Range("H138").Value = 39997.8
$ws.Range("J138").Value = 39997.8
$ws.Range("L138").Value = 39997.8
$ws.Range("N138").Value = -50277.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 199309.33
$ws.Range("J28").Value = 199309.33
$ws.Range("L28").Value = 199309.33
$ws.Range("N28").Value = -199799.33
$ws.Range("H31").Value = 3728.926
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3728.926
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3728.926
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4318.925999999999
$ws.Range("H34").Value = 3728.926
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3728.926
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3728.926
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -4132.925999999999
$ws.Range("H58").Value = 1466.9354
$ws.Range("I58").Value = 1428.1666
$ws.Range("K58").Value = 1428.1666
$ws.Range("M58").Value = -1225.1666
$ws.Range("H132").Value = 2651.743
$ws.Range("I132").Value = 2337.862
$ws.Range("K132").Value = 7013.586
$ws.Range("M132").Value = -4483.586
$ws.Range("H134").Value = 394775.75
$ws.Range("I134").Value = 509311.66
$ws.Range("J134").Value = 2081.1428
$ws.Range("K134").Value = 1527934.98
$ws.Range("L134").Value = 6243.428400000001
$ws.Range("M134").Value = -1525399.98
$ws.Range("N134").Value = -11313.4284
$ws.Range("H135").Value = 94587.22
$ws.Range("J135").Value = 94587.22
$ws.Range("L135").Value = 94587.22
$ws.Range("N135").Value = -104727.22
$ws.Range("H136").Value = 1466.9354
$ws.Range("I136").Value = 1428.1666
$ws.Range("K136").Value = 4284.4998
$ws.Range("M136").Value = -1734.4998
$ws.Range("H138").Value = 53990.816
$ws.Range("J138").Value = 53990.816
$ws.Range("L138").Value = 53990.816
$ws.Range("N138").Value = -64270.816
$ws.Range("H140").Value = 80754.914
$ws.Range("J140").Value = 80754.914
$ws.Range("L140").Value = 80754.914
$ws.Range("N140").Value = -91114.914

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 801.86664
$ws.Range("J113").Value = 819.8570999999999
$ws.Range("L113").Value = 2459.5713
$ws.Range("N113").Value = -6799.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 14508.85
$ws.Range("J123").Value = 14508.85
$ws.Range("L123").Value = 14508.85
$ws.Range("N123").Value = -19408.85
$ws.Range("H132").Value = 2292.1082
$ws.Range("I132").Value = 2076.087
$ws.Range("J132").Value = 2647
$ws.Range("K132").Value = 6228.261
$ws.Range("L132").Value = 7941
$ws.Range("M132").Value = -3698.261
$ws.Range("N132").Value = -13001
$ws.Range("H135").Value = 48575.3
$ws.Range("J135").Value = 48575.3
$ws.Range("L135").Value = 48575.3
$ws.Range("N135").Value = -58715.3
$ws.Range("H140").Value = 50666.25
$ws.Range("J140").Value = 50666.25
$ws.Range("L140").Value = 50666.25
$ws.Range("N140").Value = -61026.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 15471.429
$ws.Range("I5").Value = 16250
$ws.Range("J5").Value = 14433.333
$ws.Range("K5").Value = 16250
$ws.Range("L5").Value = 14433.333
$ws.Range("M5").Value = -16137
$ws.Range("N5").Value = -14659.333
$ws.Range("H93").Value = 1000.25
$ws.Range("I93").Value = 1000.25
$ws.Range("K93").Value = 1000.25
$ws.Range("M93").Value = 247.75
$ws.Range("H132").Value = 5870.037
$ws.Range("J132").Value = 4122.375
$ws.Range("L132").Value = 12367.125
$ws.Range("N132").Value = -17427.125
$ws.Range("H135").Value = 58211.332
$ws.Range("J135").Value = 58211.332
$ws.Range("L135").Value = 58211.332
$ws.Range("N135").Value = -68351.33199999999
$ws.Range("H136").Value = 3800.8096
$ws.Range("I136").Value = 3884.3948
$ws.Range("K136").Value = 11653.1844
$ws.Range("M136").Value = -9103.1844
$ws.Range("H137").Value = 85060
$ws.Range("J137").Value = 85060
$ws.Range("L137").Value = 85060
$ws.Range("N137").Value = -95260
$ws.Range("H139").Value = 37704.855
$ws.Range("J139").Value = 37704.855
$ws.Range("L139").Value = 37704.855
$ws.Range("N139").Value = -47984.855
$ws.Range("H141").Value = 44543
$ws.Range("J141").Value = 44543
$ws.Range("L141").Value = 44543
$ws.Range("N141").Value = -54903

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3394
$ws.Range("I132").Value = 2787.5557
$ws.Range("J132").Value = 4000.4443
$ws.Range("K132").Value = 8362.667099999999
$ws.Range("L132").Value = 12001.3329
$ws.Range("M132").Value = -5832.667099999999
$ws.Range("N132").Value = -17061.3329
$ws.Range("H133").Value = 50478.668
$ws.Range("J133").Value = 50478.668
$ws.Range("L133").Value = 50478.668
$ws.Range("N133").Value = -60598.668
$ws.Range("H136").Value = 1470.1637
$ws.Range("I136").Value = 1494.85
$ws.Range("J136").Value = 1404.3334
$ws.Range("K136").Value = 4484.549999999999
$ws.Range("L136").Value = 4213.0002
$ws.Range("M136").Value = -1934.549999999999
$ws.Range("N136").Value = -9313.0002
$ws.Range("H139").Value = 57456.332
$ws.Range("J139").Value = 57456.332
$ws.Range("L139").Value = 57456.332
$ws.Range("N139").Value = -67736.33199999999
$ws.Range("H141").Value = 76751.414
$ws.Range("J141").Value = 76751.414
$ws.Range("L141").Value = 76751.414
$ws.Range("N141").Value = -87111.414
